$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.760.81"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "3.512.91"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.665"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "4.029.91"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "621.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.01%  "
$ws.Range("D16").Value = "69.904.34"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "3.494.13"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "109.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  +6.33%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "524.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("D37").Value = "3.669.67"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +5.85%  "
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").Value = "0.0₃0773"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0469"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.143"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000241"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.53%  "
